$d = $word.ActiveDocument

# Locate the paragraph that contains "Journal of Clinical Oncology" (item 28).
$found = $d.Content.Find.Execute("Journal of Clinical Oncology", $true, $false, $false,
                                  $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Journal of Clinical Oncology' paragraph"
}

# Resolve the 1-based paragraph index by scanning the Paragraphs collection.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Journal of Clinical Oncology*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph index"
}

# Insert a new paragraph right after it.
$rng = $d.Paragraphs($targetIndex).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# Fill the new paragraph with "29. " + "Aging and Mental Health" (two runs).
$p1 = $d.Paragraphs($targetIndex + 1)
$r1 = $p1.Range
$r1.End = $r1.End - 1
$r1.Text = "29. "

$r1b = $d.Paragraphs($targetIndex + 1).Range
$r1b.End = $r1b.End - 1
$r1b.Collapse(0)
$r1b.InsertAfter("Aging and Mental Health")

# Insert a second new (blank) paragraph right after the "29." paragraph,
# inheriting the same indentation (firstLine=720) formatting.
$r2 = $d.Paragraphs($targetIndex + 1).Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# The freshly inserted blank paragraph carries a stray empty run; clear it
# so the paragraph matches the plain "<w:p><w:pPr>...</w:pPr></w:p>" shape
# used by the other blank paragraphs already in the document.
$p2 = $d.Paragraphs($targetIndex + 2)
$r3 = $p2.Range
$r3.End = $r3.End - 1
$r3.Text = "x"
$r3b = $d.Paragraphs($targetIndex + 2).Range
$r3b.End = $r3b.End - 1
$r3b.Text = ""

Write-Host "Inserted '29. Aging and Mental Health' after paragraph $targetIndex."
